$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.233.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.027.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.12%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.87'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.16'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.68%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.024.44'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.66'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.38%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.21%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.24%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.96%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.296.67'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.532.01'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.22'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.83%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.46'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +19.33%  '

$ws.Range("B20").Value = 'WrappedEther'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.033.62'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '468.44'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.67%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.37'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.02'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.73'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.03'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.22%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.19'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.43'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.23%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.21%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +7.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0993'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.20'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.11%  '

$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.992'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.86'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.15'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +9.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.05'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.46'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.311'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.69%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.89%  '

$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.63'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.47%  '

$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.84'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0360'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '378.34'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.704.15'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.37'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.44'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.23'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.09%  '
